$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4 currently holds phone number 71652621 as text (inline string); convert to numeric
$ws.Range("A4").Value = 71652621

# Add new row 5 for the new payment.
# A5 is the phone number stored as text (matches A1:A3 source data style),
# so prefix with an apostrophe to force text entry instead of a number,
# then restore the default style so only the value type changes.
$ws.Range("A5").Value = "'71652621"
$ws.Range("A5").Style = "Normal"
$ws.Range("B5").Value = 71
$ws.Range("C5").Value = "Cash"
$ws.Range("D5").Value = "2025-08-15T09:33:54"
